$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.670.97'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.745.48'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -5.64%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '236.70'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -9.65%  '
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4905'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -8.62%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.56'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -7.82%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2485'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -22.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.05967'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -15.53%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.745.30'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.72%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.06789'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -13.22%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.74'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -22.71%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.456'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -11.86%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.10'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -13.87%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.5608'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -27.72%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.001'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('E18').Value = '  +0.04%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '25.722.92'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.61%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.40'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -19.57%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.000006532'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -18.61%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.966.05'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -6.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.970'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -14.54%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.004'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -17.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.843'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -16.66%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '136.32'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.56%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.485'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -12.52%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.803'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -18.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '14.63'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -14.79%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '101.68'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -9.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.749'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -12.70%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08007'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -8.62%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.301'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -19.78%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04383'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -10.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.000'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.584'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -10.88%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9806'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -14.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.6046'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -17.99%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.676'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -14.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.999'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -14.87%  '
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('E42').Value = '  -5.62%  '
$ws.Range('E43').Value = '  -14.70%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.7556'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -17.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.139'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -13.31%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.3692'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -23.74%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05109'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -12.49%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1068'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -14.73%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '29.99'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -14.63%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '52.36'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -13.68%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.851'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -24.40%  '
